$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Week 6" header now has real dates, and the column carries real data
# (was a zero-filled placeholder column before the week happened).
$ws.Range("G1").Value = "Week 6 -- June 6 - 12"
$ws.Range("G2").Value = 455.4
$ws.Range("G3").Value = 238.6
$ws.Range("G4").Value = 337.2
$ws.Range("G5").Value = 43.3
$ws.Range("G6").Value = 388.1
$ws.Range("G7").Value = 243.2
$ws.Range("G8").Value = 169.7
$ws.Range("G10").Value = 41
$ws.Range("G11").Value = 373.4

# Match the "Week N -- dates" header column width/best-fit look.
$ws.Columns.Item(7).ColumnWidth = 17.3496875

# Blank cells (number-formatted, no value) scattered below the table --
# leftover formatting from selections/drag made while reviewing new rows.
foreach ($r in 13,14,15,16,17,18,21,22,27,29,31,33) {
    $ws.Range("B$r").NumberFormat = "0.0"
}

# Move the active selection, drop the old scrolled-down top-left cell.
$ws.Range("G15").Select() | Out-Null
